$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit performs a cyclic rotation of the data (excluding columns that
# stay identical) across rows 4, 5 and 6:
#   new row4 <- old row6
#   new row5 <- old row4
#   new row6 <- old row5
# Columns affected: A, B, D, E, F, G, H, Q, R, S

function Get-RowData($row) {
    return @{
        A = $ws.Cells.Item($row, 1).Value2
        B = $ws.Cells.Item($row, 2).Value2
        D = $ws.Cells.Item($row, 4).Value2
        E = $ws.Cells.Item($row, 5).Value2
        F = $ws.Cells.Item($row, 6).Value2
        G = $ws.Cells.Item($row, 7).Value2
        H = $ws.Cells.Item($row, 8).Value2
        Q = $ws.Cells.Item($row, 17).Value2
        R = $ws.Cells.Item($row, 18).Value2
        S = $ws.Cells.Item($row, 19).Value2
    }
}

function Set-RowData($row, $data) {
    $ws.Cells.Item($row, 1).Value2 = $data.A
    $ws.Cells.Item($row, 2).Value2 = $data.B
    $ws.Cells.Item($row, 4).Value2 = $data.D
    $ws.Cells.Item($row, 5).Value2 = $data.E
    $ws.Cells.Item($row, 6).Value2 = $data.F
    $ws.Cells.Item($row, 7).Value2 = $data.G
    $ws.Cells.Item($row, 8).Value2 = $data.H
    $ws.Cells.Item($row, 17).Value2 = $data.Q
    $ws.Cells.Item($row, 18).Value2 = $data.R
    $ws.Cells.Item($row, 19).Value2 = $data.S
}

$row4 = Get-RowData 4
$row5 = Get-RowData 5
$row6 = Get-RowData 6

Set-RowData 4 $row6
Set-RowData 5 $row4
Set-RowData 6 $row5
